# Weekly update: insert a new price record (row) for Betarraga at the top
# of the existing date-ordered block (row 226), pushing the existing rows
# 226-262 down to 227-263.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 226; this shifts the former rows
# 226..262 down to 227..263 and extends the sheet dimension automatically.
$ws.Rows("226:226").Insert()

# Populate the newly inserted row 226 with the new weekly record.
$row = 226
$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value = 44474
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100114014
$ws.Cells.Item($row, 7).Value = "Betarraga"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 120
$ws.Cells.Item($row, 11).Value = 700
$ws.Cells.Item($row, 12).Value = 750
$ws.Cells.Item($row, 13).Value = 725
$ws.Cells.Item($row, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item($row, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value = 145
$ws.Cells.Item($row, 17).Value = 5
$ws.Cells.Item($row, 18).Value = "Hortaliza"
